$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update task statuses (bcrypt integration work):
# - "Use bcrypt for hashing and storing passwords*" (row 14) is now Done
# - "Add "Profile" page" (row 5) and "Transition backend to use database..." (row 15) are In progress
$ws.Range("D5").Value2 = "In progress"
$ws.Range("D14").Value2 = "Done"
$ws.Range("D15").Value2 = "In progress"

# Remove the now-obsolete reference row that only held the bcrypt article link
$ws.Rows.Item(31).Delete()

# Reset the view: scroll back to the top and move the selection
$ws.Range("E6").Select()
